$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6603181522184486, 1.065759637188209, 0.9112048919340322, 0.6575963718820862, 1.071598303772011, 0.9070294784580499)
    3  = @(0.3950008508465115, 0.1587301587301587, -0.4402577480568004, 0.4081632653061225, 0.1559849851160945, -0.4308390022675737)
    4  = @(-0.4767228451427073, -1.26984126984127, -2.031350183375844, -0.4761904761904762, -1.26818552774021, -2.040816326530612)
    5  = @(0.6237795558828675, 0.5895691609977325, -0.002692718193406324, 0.6122448979591837, 0.6050531937043389, 0)
    6  = @(0.7086389277650147, 1.292517006802721, 1.541215953830037, 0.7029478458049886, 1.304139074833224, 1.541950113378685)
    7  = @(0.2473850426869561, 0.06802721088435373, -0.3802581138671877, 0.2494331065759637, 0.08084678929578659, -0.3854875283446712)
    8  = @(0.659905908185257, 0.9297052154195011, 0.5145133491967829, 0.6575963718820862, 0.9324018116487143, 0.5215419501133787)
    9  = @(0.5137392823977042, 0.7256235827664399, 0.4907005132262501, 0.5215419501133787, 0.7141090169048173, 0.4988662131519275)
    10 = @(0.6261701570170065, 1.08843537414966, 1.165374245474661, 0.6349206349206349, 1.07619612623426, 1.179138321995465)
    11 = @(0.7478644345922981, 1.26984126984127, 0.8777519492556363, 0.7482993197278912, 1.273801470238159, 0.8843537414965986)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
